$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A56").Value = "Nivel socioeconómico"
$ws.Range("B56").Value = "nivel_socioeconomico"

$ws.Range("A57").Value = "Tipo de hogar"
$ws.Range("B57").Value = "tipo_de_hogar"

$ws.Range("A58").Value = "Menores en el hogar"
$ws.Range("B58").Value = "menores_en_el_hogar"
